$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 31642.2046022216
$ws.Range("F2").Value = 30152.01213216448
$ws.Range("H2").Value = 21.2336360825986
$ws.Range("I2").Value = 474634.3916791933
$ws.Range("J2").Value = 318.5054288061194

# Row 3
$ws.Range("D3").Value = 31642.2046022216
$ws.Range("F3").Value = 30152.01213216448
$ws.Range("H3").Value = 21.2336360825986
$ws.Range("I3").Value = 471032.0562202815
$ws.Range("J3").Value = 316.088066263163

# Row 4
$ws.Range("D4").Value = 31642.2046022216
$ws.Range("F4").Value = 30152.01213216448
$ws.Range("H4").Value = 21.2336360825986
$ws.Range("I4").Value = 471032.0562202815
$ws.Range("J4").Value = 316.088066263163
